$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text updates
$ws.Range("B2").Value = "<both>"
$ws.Range("B3").Value = "<contro>"
$ws.Range("B5").Value = "<yil>"
$ws.Range("B6").Value = "<now>"
$ws.Range("B7").Value = "<life>"
$ws.Range("B11").Value = "<been>"
$ws.Range("B14").Value = "<they>"
$ws.Range("B15").Value = "<in>"
$ws.Range("B17").Value = "<sherrom>"

# Column C value updates
$ws.Range("C2").Value = 34
$ws.Range("C4").Value = 38
$ws.Range("C5").Value = 32
$ws.Range("C6").Value = 33
$ws.Range("C7").Value = 34
$ws.Range("C8").Value = 35
$ws.Range("C9").Value = 38
$ws.Range("C11").Value = 30
$ws.Range("C12").Value = 26
$ws.Range("C13").Value = 36
$ws.Range("C14").Value = 39
$ws.Range("C15").Value = 32
$ws.Range("C16").Value = 37
$ws.Range("C17").Value = 36
$ws.Range("C18").Value = 29
